$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the new "2 conductor" cable data
$ws.Range("B2").Value = "100+00"
$ws.Range("C2").Value = "200+00"
$ws.Range("E2").Value = "2C#2"
$ws.Range("G2").Value = 2

# H2/J2 hold percentage-looking values that must remain plain text (like the
# rest of the sheet) rather than being auto-converted to numeric percentages.
# Enter them as formulas producing the literal string, then collapse the
# formula down to a static value via copy/paste-values so the cell ends up
# as plain text without disturbing its existing style.
$ws.Range("H2").Formula = '="38.6%"'
$ws.Range("H2").Copy()
$ws.Range("H2").PasteSpecial(-4163)

$ws.Range("I2").Value = 2.5

$ws.Range("J2").Formula = '="24.70%"'
$ws.Range("J2").Copy()
$ws.Range("J2").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# Remove rows 3 through 7 (these held the merged/continuation rows for the
# old multi-pull cable table, which is no longer needed)
$ws.Range("A3:J7").EntireRow.Delete()
